$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CAN_BC")

# Row 41: keep F/H, clear G and I
$ws.Cells.Item(41, 6).Value = "11/3/1"
$ws.Cells.Item(41, 7).ClearContents()
$ws.Cells.Item(41, 8).Value = "Ref.T.HK1"
$ws.Cells.Item(41, 9).ClearContents()

# Row 42: shift up - was placeholder row (F="???"), now gets the "11/3/2" address and "Ref.T.HK3," remark
$ws.Cells.Item(42, 6).Value = "11/3/2"
$ws.Cells.Item(42, 8).Value = "Ref.T.HK3,"

# Row 43
$ws.Cells.Item(43, 6).Value = "11/3/3"
$ws.Cells.Item(43, 7).ClearContents()
$ws.Cells.Item(43, 8).Value = "Ref.T.HK4,"

# Row 44
$ws.Cells.Item(44, 6).Value = "11/3/4"
$ws.Cells.Item(44, 7).ClearContents()
$ws.Cells.Item(44, 8).Value = "Ex.Sw.HK1"

# Row 45
$ws.Cells.Item(45, 6).Value = "11/3/5"
$ws.Cells.Item(45, 8).Value = "Ex.Sw.HK3"

# Row 46
$ws.Cells.Item(46, 6).Value = "11/3/6"
$ws.Cells.Item(46, 8).Value = "Ex.Sw.HK4"

# Row 47: Ex.Sw.HK5 removed
$ws.Cells.Item(47, 6).Value = "11/3/7"
$ws.Cells.Item(47, 8).ClearContents()

# Row 48: Ex.Sw.HK4 removed (duplicate, now only exists at row46)
$ws.Cells.Item(48, 6).Value = "11/3/8"
$ws.Cells.Item(48, 8).ClearContents()

# Row 49: Ex.Sw.HK3 removed
$ws.Cells.Item(49, 6).Value = "11/3/9"
$ws.Cells.Item(49, 8).ClearContents()

# Row 50 unchanged
$ws.Cells.Item(50, 6).Value = "11/3/10"
